$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ofertas")
$rng = $ws.Range("E2:F2")
$rng.NumberFormat = "mm-dd-yy"
$rng.Value = 45437
Write-Output $ws.Range("E2").Text
Write-Output $ws.Range("F2").Text
